# Apply crypto price/volume updates from the commit diff.
# Use NumberFormat "@" to force text, then restore "Normal" style so
# no residual style index is left behind on cells that had none.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "56.608.36"
Set-TextValue $ws.Range("E2") "  +0.04%  "
Set-TextValue $ws.Range("D3") "2.312.46"
Set-TextValue $ws.Range("E3") "  -0.44%  "
Set-TextValue $ws.Range("D5") "519.18"
Set-TextValue $ws.Range("E5") "  +0.64%  "
Set-TextValue $ws.Range("D6") "132.07"
Set-TextValue $ws.Range("E6") "  -2.48%  "
Set-TextValue $ws.Range("D7") "0.995"
Set-TextValue $ws.Range("E7") "  -0.12%  "
Set-TextValue $ws.Range("D8") "0.533"
Set-TextValue $ws.Range("E8") "  -0.77%  "
Set-TextValue $ws.Range("D9") "2.332.19"
Set-TextValue $ws.Range("E9") "  -0.49%  "
Set-TextValue $ws.Range("D10") "0.100"
Set-TextValue $ws.Range("E10") "  -1.87%  "
Set-TextValue $ws.Range("E11") "  +0.03%  "
Set-TextValue $ws.Range("D12") "5.26"
Set-TextValue $ws.Range("E12") "  -1.73%  "
Set-TextValue $ws.Range("D13") "0.339"
Set-TextValue $ws.Range("E13") "  -1.37%  "
Set-TextValue $ws.Range("B14") "Avalanche"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D14") "23.43"
Set-TextValue $ws.Range("E14") "  -2.04%  "
Set-TextValue $ws.Range("B15") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "2.723.04"
Set-TextValue $ws.Range("E15") "  -0.61%  "
Set-TextValue $ws.Range("D16") "56.543.60"
Set-TextValue $ws.Range("E16") "  -0.18%  "
Set-TextValue $ws.Range("E17") "  -1.64%  "
Set-TextValue $ws.Range("D18") "2.327.91"
Set-TextValue $ws.Range("E18") "  +0.35%  "
Set-TextValue $ws.Range("D19") "334.11"
Set-TextValue $ws.Range("E19") "  +2.42%  "
Set-TextValue $ws.Range("D20") "10.37"
Set-TextValue $ws.Range("E20") "  -1.44%  "
Set-TextValue $ws.Range("D21") "4.15"
Set-TextValue $ws.Range("E21") "  -1.75%  "
Set-TextValue $ws.Range("D22") "6.74"
Set-TextValue $ws.Range("E22") "  +2.15%  "
Set-TextValue $ws.Range("E23") "  -0.13%  "
Set-TextValue $ws.Range("D24") "61.24"
Set-TextValue $ws.Range("E24") "  +0.75%  "
Set-TextValue $ws.Range("D25") "8.67"
Set-TextValue $ws.Range("E25") "  +8.53%  "
Set-TextValue $ws.Range("D26") "0.165"
Set-TextValue $ws.Range("E26") "  -0.22%  "
Set-TextValue $ws.Range("D27") "0.992"
Set-TextValue $ws.Range("E27") "  -0.22%  "
Set-TextValue $ws.Range("D28") "1.34"
Set-TextValue $ws.Range("E28") "  +4.39%  "
Set-TextValue $ws.Range("D29") "169.67"
Set-TextValue $ws.Range("E29") "  +0.21%  "
Set-TextValue $ws.Range("E30") "  -0.06%  "
Set-TextValue $ws.Range("D31") "0.0₃0718"
Set-TextValue $ws.Range("E31") "  -3.36%  "
Set-TextValue $ws.Range("D32") "6.12"
Set-TextValue $ws.Range("E32") "  -1.19%  "
Set-TextValue $ws.Range("D33") "18.35"
Set-TextValue $ws.Range("E33") "  -0.71%  "
Set-TextValue $ws.Range("D34") "0.998"
Set-TextValue $ws.Range("E34") "  -0.03%  "
Set-TextValue $ws.Range("D35") "0.994"
Set-TextValue $ws.Range("E35") "  -0.07%  "
Set-TextValue $ws.Range("D36") "1.25"
Set-TextValue $ws.Range("E36") "  -0.79%  "
Set-TextValue $ws.Range("D37") "3.93"
Set-TextValue $ws.Range("E37") "  -1.77%  "
Set-TextValue $ws.Range("D38") "0.884"
Set-TextValue $ws.Range("E38") "  -3.50%  "
Set-TextValue $ws.Range("D39") "1.58"
Set-TextValue $ws.Range("E39") "  +1.26%  "
Set-TextValue $ws.Range("D40") "38.79"
Set-TextValue $ws.Range("E40") "  +1.17%  "
Set-TextValue $ws.Range("D41") "147.95"
Set-TextValue $ws.Range("E41") "  +4.33%  "
Set-TextValue $ws.Range("D42") "0.375"
Set-TextValue $ws.Range("E42") "  -1.66%  "
Set-TextValue $ws.Range("D43") "287.67"
Set-TextValue $ws.Range("E43") "  +3.38%  "
Set-TextValue $ws.Range("D44") "3.58"
Set-TextValue $ws.Range("E44") "  -0.73%  "
Set-TextValue $ws.Range("D45") "5.10"
Set-TextValue $ws.Range("E45") "  -1.36%  "
Set-TextValue $ws.Range("D46") "0.0927"
Set-TextValue $ws.Range("E46") "  -1.04%  "
Set-TextValue $ws.Range("D47") "0.0500"
Set-TextValue $ws.Range("E47") "  -1.34%  "
Set-TextValue $ws.Range("D48") "0.558"
Set-TextValue $ws.Range("E48") "  -0.91%  "
Set-TextValue $ws.Range("D49") "18.37"
Set-TextValue $ws.Range("E49") "  +2.92%  "
Set-TextValue $ws.Range("D50") "0.0214"
Set-TextValue $ws.Range("E50") "  -1.94%  "
Set-TextValue $ws.Range("D51") "0.377"
Set-TextValue $ws.Range("E51") "  -0.68%  "
